$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (ds / order 2, "Communicating and transforming data") now reflects the
# re-launched course offering: year/crn become combined ranges, and the
# website + repo point at the new "dataviz-2020" materials instead of the
# old "communicate_transform_data" ones.
$ws.Range("D3").Value = "2019/2020"
$ws.Range("E3").Value = "27553/27120"
$ws.Range("H3").Value = "https://uo-datasci-specialization.github.io/c2-dataviz-2020/"
$ws.Range("I3").Value = "https://github.com/uo-datasci-specialization/c2-dataviz-2020"

# Move the active selection to G13, matching the saved cursor position.
$ws.Range("G13").Select()
